# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Femacal de La Calera - Frutilla" just above
# the existing row 146, shifting the previous rows 146-150 down to 148-152.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 146.
$ws.Rows("146:147").Insert()

$rowDate = [DateTime]"2021-11-09"

# --- New row 146 -------------------------------------------------------
$ws.Range("A146").Value = 3
$ws.Range("B146").Value = "Femacal de La Calera"
$ws.Range("C146").Value = "Coquimbo"
$ws.Range("D146").Value = $rowDate
$ws.Range("E146").Value = 5
$ws.Range("F146").Value = "Fruta"
$ws.Range("G146").Value = 100101
$ws.Range("H146").Value = "Berries"
$ws.Range("I146").Value = 100112025
$ws.Range("J146").Value = "Frutilla"
$ws.Range("K146").Value = "Sin especificar"
$ws.Range("L146").Value = "Especial"
$ws.Range("M146").Value = 125
$ws.Range("N146").Value = 5500
$ws.Range("O146").Value = 6000
$ws.Range("P146").Value = 5760
$ws.Range("Q146").Value = "$/bandeja 7 kilos"
$ws.Range("R146").Value = "Provincia de Melipilla"
$ws.Range("S146").Value = 823
$ws.Range("T146").Value = 7

# --- New row 147 -------------------------------------------------------
$ws.Range("A147").Value = 3
$ws.Range("B147").Value = "Femacal de La Calera"
$ws.Range("C147").Value = "Coquimbo"
$ws.Range("D147").Value = $rowDate
$ws.Range("E147").Value = 5
$ws.Range("F147").Value = "Fruta"
$ws.Range("G147").Value = 100101
$ws.Range("H147").Value = "Berries"
$ws.Range("I147").Value = 100112025
$ws.Range("J147").Value = "Frutilla"
$ws.Range("K147").Value = "Sin especificar"
$ws.Range("L147").Value = "Segunda"
$ws.Range("M147").Value = 50
$ws.Range("N147").Value = 4000
$ws.Range("O147").Value = 4000
$ws.Range("P147").Value = 4000
$ws.Range("Q147").Value = "$/bandeja 7 kilos"
$ws.Range("R147").Value = "Provincia de Melipilla"
$ws.Range("S147").Value = 571
$ws.Range("T147").Value = 7
